$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "67.178.41"
$ws.Cells.Item(2, 5).Value = "  +0.29%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "3.108.07"
$ws.Cells.Item(3, 5).Value = "  -0.22%  "
$ws.Cells.Item(4, 5).Value = "  +0.02%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "579.05"
$ws.Cells.Item(5, 5).Value = "  -0.08%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "173.59"
$ws.Cells.Item(6, 5).Value = "  +0.33%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "1.00"
$ws.Cells.Item(7, 5).Value = "  +0.09%  "
$ws.Cells.Item(8, 5).Value = "  -0.80%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "6.53"
$ws.Cells.Item(9, 5).Value = "  +1.14%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.153"
$ws.Cells.Item(10, 5).Value = "  -1.79%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.477"
$ws.Cells.Item(11, 5).Value = "  -1.16%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.0000247"
$ws.Cells.Item(12, 5).Value = "  -0.97%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "36.73"
$ws.Cells.Item(13, 5).Value = "  -1.57%  "
$ws.Cells.Item(14, 5).Value = "  -1.69%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "3.623.38"
$ws.Cells.Item(15, 5).Value = "  -0.10%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "67.119.01"
$ws.Cells.Item(16, 5).Value = "  +0.23%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "7.08"
$ws.Cells.Item(17, 5).Value = "  -1.66%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "3.113.38"
$ws.Cells.Item(18, 5).Value = "  +0.11%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "16.56"
$ws.Cells.Item(19, 5).Value = "  +1.80%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "491.33"
$ws.Cells.Item(20, 5).Value = "  +0.77%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "0.702"
$ws.Cells.Item(21, 5).Value = "  -2.44%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "7.84"
$ws.Cells.Item(22, 5).Value = "  +3.44%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "83.94"
$ws.Cells.Item(23, 5).Value = "  -0.75%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "13.07"
$ws.Cells.Item(24, 5).Value = "  -2.29%  "
$ws.Cells.Item(25, 5).Value = "  -3.30%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "10.58"
$ws.Cells.Item(26, 5).Value = "  +4.88%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "7.91"
$ws.Cells.Item(28, 5).Value = "  -1.85%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "2.35"
$ws.Cells.Item(29, 5).Value = "  -3.27%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "2.66"
$ws.Cells.Item(30, 5).Value = "  -1.03%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "28.32"
$ws.Cells.Item(31, 5).Value = "  -2.48%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.113"
$ws.Cells.Item(32, 5).Value = "  -1.37%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.0₃0934"
$ws.Cells.Item(33, 5).Value = "  -7.20%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "1.00"
$ws.Cells.Item(34, 5).Value = "  +0.13%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "5.83"
$ws.Cells.Item(35, 5).Value = "  -1.97%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.970"
$ws.Cells.Item(36, 5).Value = "  -1.88%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "47.14"
$ws.Cells.Item(37, 5).Value = "  -0.74%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "2.04"
$ws.Cells.Item(38, 5).Value = "  -3.82%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.308"
$ws.Cells.Item(39, 5).Value = "  -2.58%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.123"
$ws.Cells.Item(40, 5).Value = "  +0.72%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "8.46"
$ws.Cells.Item(41, 5).Value = "  -2.69%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "388.16"
$ws.Cells.Item(42, 5).Value = "  +0.59%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "2.798.64"
$ws.Cells.Item(43, 5).Value = "  -1.69%  "
$ws.Cells.Item(44, 5).Value = "  -8.61%  "
$ws.Cells.Item(45, 5).Value = "  -2.28%  "
$ws.Cells.Item(46, 5).Value = "  -1.28%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "25.05"
$ws.Cells.Item(48, 5).Value = "  -0.52%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "2.20"
$ws.Cells.Item(49, 5).Value = "  -1.42%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.107"
$ws.Cells.Item(50, 5).Value = "  -1.55%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "6.69"
$ws.Cells.Item(51, 5).Value = "  -2.94%  "
